$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 2350
$ws.Range("I94").Value = 2350
$ws.Range("K94").Value = 2350
$ws.Range("M94").Value = -1899
$ws.Range("H98").Value = 1130.7878
$ws.Range("I98").Value = 910
$ws.Range("J98").Value = 1719.5555
$ws.Range("K98").Value = 910
$ws.Range("L98").Value = 1719.5555
$ws.Range("M98").Value = 588
$ws.Range("N98").Value = -4715.5555
$ws.Range("H122").Value = 1130.7878
$ws.Range("I122").Value = 910
$ws.Range("J122").Value = 1719.5555
$ws.Range("K122").Value = 2730
$ws.Range("L122").Value = 5158.666499999999
$ws.Range("M122").Value = -280
$ws.Range("N122").Value = -10058.6665
$ws.Range("H132").Value = 1028.4769
$ws.Range("I132").Value = 721.3774
$ws.Range("J132").Value = 2384.8333
$ws.Range("K132").Value = 2164.1322
$ws.Range("L132").Value = 7154.499899999999
$ws.Range("M132").Value = 365.8678
$ws.Range("N132").Value = -12214.4999
$ws.Range("H135").Value = 1340.3334
$ws.Range("I135").Value = 1025.96
$ws.Range("K135").Value = 9233.639999999999
$ws.Range("M135").Value = -6698.639999999999
$ws.Range("H137").Value = 1295.7322
$ws.Range("I137").Value = 1045.2325
$ws.Range("J137").Value = 2124.3076
$ws.Range("K137").Value = 3135.6975
$ws.Range("L137").Value = 6372.9228
$ws.Range("M137").Value = -585.6975000000002
$ws.Range("N137").Value = -11472.9228
$ws.Range("H138").Value = 2192.4268
$ws.Range("I138").Value = 870.5333000000001
$ws.Range("J138").Value = 3800.1353
$ws.Range("K138").Value = 2611.5999
$ws.Range("L138").Value = 11400.4059
$ws.Range("M138").Value = 2528.4001
$ws.Range("N138").Value = -21680.4059
$ws.Range("H141").Value = 1461.74
$ws.Range("I141").Value = 968.89746
$ws.Range("K141").Value = 2906.69238
$ws.Range("M141").Value = 2273.30762

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3719.08
$ws.Range("I32").Value = 2922.15
$ws.Range("J32").Value = 6906.8
$ws.Range("K32").Value = 2922.15
$ws.Range("L32").Value = 6906.8
$ws.Range("M32").Value = -2635.15
$ws.Range("N32").Value = -7480.8
$ws.Range("H97").Value = 1007.5238
$ws.Range("I97").Value = 872.3333
$ws.Range("J97").Value = 1345.5
$ws.Range("K97").Value = 872.3333
$ws.Range("L97").Value = 1345.5
$ws.Range("M97").Value = -376.3333
$ws.Range("N97").Value = -2337.5
$ws.Range("H122").Value = 2565219.5
$ws.Range("I122").Value = 2849688.2
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 8549064.600000001
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -8546614.600000001
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 1668903.4
$ws.Range("I132").Value = 1366.2916
$ws.Range("J132").Value = 8339051.5
$ws.Range("K132").Value = 4098.8748
$ws.Range("L132").Value = 25017154.5
$ws.Range("M132").Value = -1568.8748
$ws.Range("N132").Value = -25022214.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1258.6364
$ws.Range("I94").Value = 439.23077
$ws.Range("K94").Value = 439.23077
$ws.Range("M94").Value = 11.76922999999999
$ws.Range("H134").Value = 1554.277
$ws.Range("I134").Value = 1176.4651
$ws.Range("J134").Value = 2292.7273
$ws.Range("K134").Value = 3529.3953
$ws.Range("L134").Value = 6878.1819
$ws.Range("M134").Value = -994.3952999999997
$ws.Range("N134").Value = -11948.1819

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6645.5474
$ws.Range("I31").Value = 1699.8055
$ws.Range("K31").Value = 1699.8055
$ws.Range("M31").Value = -1404.8055
$ws.Range("H34").Value = 6645.5474
$ws.Range("I34").Value = 1699.8055
$ws.Range("K34").Value = 1699.8055
$ws.Range("M34").Value = -1497.8055
$ws.Range("H69").Value = 5776.4
$ws.Range("I69").Value = 1470.5
$ws.Range("J69").Value = 23000
$ws.Range("K69").Value = 1470.5
$ws.Range("L69").Value = 23000
$ws.Range("M69").Value = -721.5
$ws.Range("N69").Value = -24498
$ws.Range("H72").Value = 5776.4
$ws.Range("I72").Value = 1470.5
$ws.Range("J72").Value = 23000
$ws.Range("K72").Value = 4411.5
$ws.Range("L72").Value = 69000
$ws.Range("M72").Value = -667.5
$ws.Range("N72").Value = -76488
$ws.Range("H99").Value = 5216070.5
$ws.Range("I99").Value = 7348.8335
$ws.Range("J99").Value = 20842236
$ws.Range("K99").Value = 7348.8335
$ws.Range("L99").Value = 20842236
$ws.Range("M99").Value = -5850.8335
$ws.Range("N99").Value = -20845232
$ws.Range("H105").Value = 2061.2083
$ws.Range("I105").Value = 2264.647
$ws.Range("J105").Value = 1567.1428
$ws.Range("K105").Value = 2264.647
$ws.Range("L105").Value = 1567.1428
$ws.Range("M105").Value = -517.6469999999999
$ws.Range("N105").Value = -5061.1428
$ws.Range("H122").Value = 1517.5883
$ws.Range("I122").Value = 1032.3334
$ws.Range("J122").Value = 2063.5
$ws.Range("K122").Value = 3097.0002
$ws.Range("L122").Value = 6190.5
$ws.Range("M122").Value = -647.0001999999999
$ws.Range("N122").Value = -11090.5
$ws.Range("H126").Value = 5216070.5
$ws.Range("I126").Value = 7348.8335
$ws.Range("J126").Value = 20842236
$ws.Range("K126").Value = 22046.5005
$ws.Range("L126").Value = 62526708
$ws.Range("M126").Value = -19576.5005
$ws.Range("N126").Value = -62531648
$ws.Range("H132").Value = 2326.75
$ws.Range("I132").Value = 1830.2858
$ws.Range("J132").Value = 3195.5625
$ws.Range("K132").Value = 5490.857400000001
$ws.Range("L132").Value = 9586.6875
$ws.Range("M132").Value = -2960.857400000001
$ws.Range("N132").Value = -14646.6875

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 10110.417
$ws.Range("J57").Value = 10110.417
$ws.Range("L57").Value = 10110.417
$ws.Range("N57").Value = -11750.417
$ws.Range("H70").Value = 6353.9473
$ws.Range("I70").Value = 6513.5
$ws.Range("J70").Value = 5503
$ws.Range("K70").Value = 6513.5
$ws.Range("L70").Value = 5503
$ws.Range("M70").Value = -6243.5
$ws.Range("N70").Value = -6043
$ws.Range("H73").Value = 6353.9473
$ws.Range("I73").Value = 6513.5
$ws.Range("J73").Value = 5503
$ws.Range("K73").Value = 6513.5
$ws.Range("L73").Value = 5503
$ws.Range("M73").Value = -5577.5
$ws.Range("N73").Value = -7375
$ws.Range("H80").Value = 3016.6365
$ws.Range("I80").Value = 2276.6
$ws.Range("J80").Value = 3633.3333
$ws.Range("K80").Value = 2276.6
$ws.Range("L80").Value = 3633.3333
$ws.Range("M80").Value = -1278.6
$ws.Range("N80").Value = -5629.3333
$ws.Range("H83").Value = 3016.6365
$ws.Range("I83").Value = 2276.6
$ws.Range("J83").Value = 3633.3333
$ws.Range("K83").Value = 11383
$ws.Range("L83").Value = 18166.6665
$ws.Range("M83").Value = -6391
$ws.Range("N83").Value = -28150.6665
$ws.Range("H102").Value = 2104.75
$ws.Range("I102").Value = 1740.5
$ws.Range("J102").Value = 2286.875
$ws.Range("K102").Value = 1740.5
$ws.Range("L102").Value = 2286.875
$ws.Range("M102").Value = -118.5
$ws.Range("N102").Value = -5530.875
$ws.Range("H122").Value = 81911096
$ws.Range("I122").Value = 118314770
$ws.Range("J122").Value = 2822
$ws.Range("K122").Value = 354944310
$ws.Range("L122").Value = 8466
$ws.Range("M122").Value = -354941860
$ws.Range("N122").Value = -13366
$ws.Range("H123").Value = 21294.926
$ws.Range("J123").Value = 21294.926
$ws.Range("L123").Value = 21294.926
$ws.Range("N123").Value = -26194.926
$ws.Range("H126").Value = 4550.8184
$ws.Range("I126").Value = 8454.134
$ws.Range("K126").Value = 25362.402
$ws.Range("M126").Value = -22892.402
$ws.Range("H132").Value = 1815.3043
$ws.Range("I132").Value = 1374.1154
$ws.Range("J132").Value = 2388.85
$ws.Range("K132").Value = 4122.3462
$ws.Range("L132").Value = 7166.549999999999
$ws.Range("M132").Value = -1592.3462
$ws.Range("N132").Value = -12226.55

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 85525.336
$ws.Range("I7").Value = 92891.27
$ws.Range("K7").Value = 92891.27
$ws.Range("M7").Value = -92779.27
$ws.Range("H40").Value = 62504388
$ws.Range("I40").Value = 76925784
$ws.Range("J40").Value = 11665
$ws.Range("K40").Value = 76925784
$ws.Range("L40").Value = 11665
$ws.Range("M40").Value = -76925648
$ws.Range("N40").Value = -11937
$ws.Range("H100").Value = 2001.3334
$ws.Range("I100").Value = 2000
$ws.Range("J100").Value = 2004
$ws.Range("K100").Value = 2000
$ws.Range("L100").Value = 2004
$ws.Range("M100").Value = -1459
$ws.Range("N100").Value = -3086
$ws.Range("H122").Value = 2091814.8
$ws.Range("I122").Value = 2555456.2
$ws.Range("J122").Value = 911636.2
$ws.Range("K122").Value = 7666368.600000001
$ws.Range("L122").Value = 2734908.6
$ws.Range("M122").Value = -7663918.600000001
$ws.Range("N122").Value = -2739808.6
$ws.Range("H126").Value = 85525.336
$ws.Range("I126").Value = 92891.27
$ws.Range("K126").Value = 278673.81
$ws.Range("M126").Value = -276203.81
$ws.Range("H132").Value = 10084191
$ws.Range("I132").Value = 12724320
$ws.Range("J132").Value = 3699.7273
$ws.Range("K132").Value = 38172960
$ws.Range("L132").Value = 11099.1819
$ws.Range("M132").Value = -38170430
$ws.Range("N132").Value = -16159.1819
$ws.Range("H136").Value = 4068.0908
$ws.Range("I136").Value = 1665.2963
$ws.Range("K136").Value = 4995.8889
$ws.Range("M136").Value = -2445.8889

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1337.2273
$ws.Range("I113").Value = 1250.4
$ws.Range("J113").Value = 1523.2858
$ws.Range("K113").Value = 3751.2
$ws.Range("L113").Value = 4569.857400000001
$ws.Range("M113").Value = -1581.2
$ws.Range("N113").Value = -8909.857400000001
$ws.Range("H122").Value = 2212.32
$ws.Range("I122").Value = 2195.158
$ws.Range("J122").Value = 2266.6667
$ws.Range("K122").Value = 6585.474
$ws.Range("L122").Value = 6800.000100000001
$ws.Range("M122").Value = -4135.474
$ws.Range("N122").Value = -11700.0001
$ws.Range("H132").Value = 22364.914
$ws.Range("I132").Value = 25154.83
$ws.Range("J132").Value = 3300.5
$ws.Range("K132").Value = 75464.49000000001
$ws.Range("L132").Value = 9901.5
$ws.Range("M132").Value = -72934.49000000001
$ws.Range("N132").Value = -14961.5

